$wb = $excel.ActiveWorkbook

# Sheet5: add a new leaderboard row (row 6) for player "r"
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Range("A6").Value = "r"
$ws5.Range("B6").Value = 2084.0

# Sheet3: add three new leaderboard rows (rows 4-6) for player "k"
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A4").Value = "k"
$ws3.Range("B4").Value = 0.0
$ws3.Range("A5").Value = "k"
$ws3.Range("B5").Value = 3126.0
$ws3.Range("A6").Value = "k"
$ws3.Range("B6").Value = 3126.0

# Sheet1: add a new leaderboard row (row 21) for player "A" with score 1042
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A21").Value = "A"
$ws1.Range("B21").Value = 1042.0
